$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3.286832544864788
$ws.Range("C2").Value = 1.655778082260271
$ws.Range("D2").Value = 0.7527432677738641
$ws.Range("E2").Value = 0.4942365360607697
$ws.Range("G2").Value = 6.189590430959694
